$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.016.89"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "3.266.44"
$ws.Range("E3").Value = "  +4.97%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "615.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.380"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.15%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "3.265.59"
$ws.Range("E10").Value = "  +5.13%  "
$ws.Range("E11").Value = "  -7.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.197"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "95.999.08"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.874.83"
$ws.Range("E16").Value = "  +5.22%  "
$ws.Range("E17").Value = "  +2.58%  "
$ws.Range("D18").Value = "3.280.90"
$ws.Range("E18").Value = "  +5.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "475.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000202"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").Value = "3.440.48"
$ws.Range("E28").Value = "  +4.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("E31").Value = "  -7.59%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.69%  "
$ws.Range("E37").Value = "  -6.82%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "487.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.443"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("E45").Value = "  -6.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.773"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.29%  "
